$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column B ("Nama"), shifting subsequent columns left.
$ws.Columns.Item(2).Delete()

# Leave the (now-shifted) column B selected, matching Excel's behavior
# after an entire-column delete.
$ws.Columns.Item(2).Select()
